# Commit: "Created setting with nine scenarios"
#
# Adds a new worksheet "nine_scenarios" right after "scenarios_base",
# populated with a 3x3x3 full-factorial scenario table (9 scenarios +
# header row), mirroring the layout of the existing "scenarios_base"
# sheet but driven by 1/9 probabilities instead of 1/27.

$wb = $excel.ActiveWorkbook

# --- 1. Select the full data range on scenarios_base (basis for the new sheet) ---
$base = $wb.Worksheets.Item("scenarios_base")
$base.Activate()
$base.Range("A1:K28").Select()

# --- 2. Insert the new worksheet right after scenarios_base ---
$ws = $wb.Worksheets.Add($null, $base)
$ws.Name = "nine_scenarios"

# --- 3. Header row ---
$ws.Range("A1").Value = "Scenario"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Probability"
$ws.Range("D1").Value = "Cost_Established"
$ws.Range("E1").Value = "Cost_Battery"
$ws.Range("F1").Value = "Cost_Hydrogen"
$ws.Range("G1").Value = "Cost_Biofuel"
$ws.Range("H1").Value = "Maturity_Established"
$ws.Range("I1").Value = "Maturity_Battery"
$ws.Range("J1").Value = "Maturity_Hydrogen"
$ws.Range("K1").Value = "Maturity_Biofuel"
$ws.Range("A1:K1").Font.Bold = $true

# --- 4. Scenario index column (A) ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8

# --- 5. Scenario name column (B) -- full factorial BBB/OOO/... codes ---
$ws.Range("B2").Value = "BBB"
$ws.Range("B3").Value = "OOO"
$ws.Range("B4").Value = "OOP"
$ws.Range("B5").Value = "OPO"
$ws.Range("B6").Value = "OPP"
$ws.Range("B7").Value = "POO"
$ws.Range("B8").Value = "POP"
$ws.Range("B9").Value = "PPO"
$ws.Range("B10").Value = "PPP"

# --- 6. Probability column (C) -- C2 standalone, C3:C10 shared formula 1/9 ---
$ws.Range("C2").Formula = "=1/9"
$ws.Range("C3:C10").Formula = "=1/9"

# --- 7. Cost_Established column (D) -- always 1 ---
$ws.Range("D2:D10").Value = 1

# --- 8. Cost_Battery / Cost_Hydrogen / Cost_Biofuel columns (E, F, G) ---
# Row 2 (BBB) stays at the base multiplier (1); all other rows reference
# the +/-25% variability toggle on the "variability" sheet.
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("E3").Formula = "=1-variability!`$A`$2"
$ws.Range("F3").Formula = "=1-variability!`$A`$2"
$ws.Range("G3").Formula = "=1-variability!`$A`$2"

$ws.Range("E4").Formula = "=1-variability!`$A`$2"
$ws.Range("F4").Formula = "=1-variability!`$A`$2"
$ws.Range("G4").Formula = "=1+variability!`$A`$2"

$ws.Range("E5").Formula = "=1-variability!`$A`$2"
$ws.Range("F5").Formula = "=1+variability!`$A`$2"
$ws.Range("G5").Formula = "=1-variability!`$A`$2"

$ws.Range("E6").Formula = "=1-variability!`$A`$2"
$ws.Range("F6").Formula = "=1+variability!`$A`$2"
$ws.Range("G6").Formula = "=1+variability!`$A`$2"

$ws.Range("E7").Formula = "=1+variability!`$A`$2"
$ws.Range("F7").Formula = "=1-variability!`$A`$2"
$ws.Range("G7").Formula = "=1-variability!`$A`$2"

$ws.Range("E8").Formula = "=1+variability!`$A`$2"
$ws.Range("F8").Formula = "=1-variability!`$A`$2"
$ws.Range("G8").Formula = "=1+variability!`$A`$2"

$ws.Range("E9").Formula = "=1+variability!`$A`$2"
$ws.Range("F9").Formula = "=1+variability!`$A`$2"
$ws.Range("G9").Formula = "=1-variability!`$A`$2"

$ws.Range("E10").Formula = "=1+variability!`$A`$2"
$ws.Range("F10").Formula = "=1+variability!`$A`$2"
$ws.Range("G10").Formula = "=1+variability!`$A`$2"

# --- 9. Maturity columns (H, I, J, K) -- base/fast/slow text codes ---
$ws.Range("H2").Value = "base"
$ws.Range("I2").Value = "base"
$ws.Range("J2").Value = "base"
$ws.Range("K2").Value = "base"

$ws.Range("H3").Value = "base"
$ws.Range("I3").Value = "fast"
$ws.Range("J3").Value = "fast"
$ws.Range("K3").Value = "fast"

$ws.Range("H4").Value = "base"
$ws.Range("I4").Value = "fast"
$ws.Range("J4").Value = "fast"
$ws.Range("K4").Value = "slow"

$ws.Range("H5").Value = "base"
$ws.Range("I5").Value = "fast"
$ws.Range("J5").Value = "slow"
$ws.Range("K5").Value = "fast"

$ws.Range("H6").Value = "base"
$ws.Range("I6").Value = "fast"
$ws.Range("J6").Value = "slow"
$ws.Range("K6").Value = "slow"

$ws.Range("H7").Value = "base"
$ws.Range("I7").Value = "slow"
$ws.Range("J7").Value = "fast"
$ws.Range("K7").Value = "fast"

$ws.Range("H8").Value = "base"
$ws.Range("I8").Value = "slow"
$ws.Range("J8").Value = "fast"
$ws.Range("K8").Value = "slow"

$ws.Range("H9").Value = "base"
$ws.Range("I9").Value = "slow"
$ws.Range("J9").Value = "slow"
$ws.Range("K9").Value = "fast"

$ws.Range("H10").Value = "base"
$ws.Range("I10").Value = "slow"
$ws.Range("J10").Value = "slow"
$ws.Range("K10").Value = "slow"

# --- 10. View state for the new sheet ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Activate()
$ws.Range("G16").Select()

# --- 11. EV_scenario loses tab focus to the new sheet; its selection moves too ---
$ev = $wb.Worksheets.Item("EV_scenario")
$ev.Activate()
$ev.Range("F10").Select()

# --- 12. Leave nine_scenarios as the active/selected sheet, as in the target ---
$ws.Activate()
$ws.Range("G16").Select()
